$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (formerly Strike#) values.
# Update per-row K values as recalculated (K instead of Strike#).
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 0
